# Apply updated odds values to the "Jogos da Semana" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("V2").Value = 1.36

# Row 3
$ws.Range("G3").Value = 3.3
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 2.3
$ws.Range("J3").Value = 4
$ws.Range("U3").Value = 3.8
$ws.Range("Y3").Value = 1.57
$ws.Range("Z3").Value = 2.25
$ws.Range("AF3").Value = 34

# Row 4
$ws.Range("S4").Value = 2.35
$ws.Range("T4").Value = 1.57

# Row 5
$ws.Range("G5").Value = 2.38
$ws.Range("I5").Value = 2.9
$ws.Range("J5").Value = 3.2
$ws.Range("L5").Value = 3.6
$ws.Range("N5").Value = 8.5
$ws.Range("Y5").Value = 1.44
$ws.Range("Z5").Value = 2.63
$ws.Range("AA5").Value = 1.83
$ws.Range("AB5").Value = 1.83
$ws.Range("AC5").Value = 7.5
$ws.Range("AD5").Value = 11
$ws.Range("AE5").Value = 10
$ws.Range("AF5").Value = 23
$ws.Range("AG5").Value = 21
$ws.Range("AK5").Value = 15
$ws.Range("AN5").Value = 13
$ws.Range("AO5").Value = 11
$ws.Range("AP5").Value = 29
$ws.Range("AQ5").Value = 26
$ws.Range("AR5").Value = 34
$ws.Range("AS5").Value = 301
